# Generate Report for Handback
# Row 7 (bdbebf19-2cd6-4b59-b192-85806429c2c2) in both the "zh-cn" and
# "de-de" sheets moves from "handed off, awaiting handback" to
# "handback received" state: the Latest Target File / Latest Handback
# File / Latest Handback DateTime / Error Detail columns get filled in.

$wb = $excel.ActiveWorkbook

$targetDisplay = "bdbebf19-2cd6-4b59-b192-85806429c2c2.md"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("J7").Value = "bdbebf19-2cd6-4b59-b192-85806429c2c2.06223f76aca9ea017ded3f93443fa081b9a07127.zh-cn.xlf"
$wsZh.Range("K7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5e8a98c040769e54ff17614f6b933adaf656d936/e2e/bdbebf19-2cd6-4b59-b192-85806429c2c2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d9876ae98bc73c1c26694a7c65b3eb1bc1a5d680/e2e/bdbebf19-2cd6-4b59-b192-85806429c2c2.md."
$wsZh.Range("P7").Value = "2016-09-05 03:01:28"

$wsZh.Range("I7").Value = $targetDisplay
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/06223f76aca9ea017ded3f93443fa081b9a07127/e2e/bdbebf19-2cd6-4b59-b192-85806429c2c2.md", "", "", $targetDisplay)

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("J7").Value = "bdbebf19-2cd6-4b59-b192-85806429c2c2.06223f76aca9ea017ded3f93443fa081b9a07127.de-de.xlf"
$wsDe.Range("K7").Value = "2016-09-05 03:01:53"
$wsDe.Range("P7").Value = "2016-09-05 03:01:28"

$wsDe.Range("I7").Value = $targetDisplay
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/06223f76aca9ea017ded3f93443fa081b9a07127/e2e/bdbebf19-2cd6-4b59-b192-85806429c2c2.md", "", "", $targetDisplay)
